$wb = $excel.ActiveWorkbook
$wsTranslations = $wb.Worksheets.Item("Translations")
$wsQuestion = $wb.Worksheets.Item("@@_question")

# Insert a new "Variable" column (B) on both sheets, shifting existing
# columns one to the right.
$wsTranslations.Columns.Item(2).Insert()
$wsQuestion.Columns.Item(2).Insert()

# Header
$wsTranslations.Range("B1").Value = "Variable"
$wsQuestion.Range("B1").Value = "Variable"

# Populate the new column with the question variable name for every data row.
$wsTranslations.Range("B2").Value = "s1"
$wsTranslations.Range("B3").Value = "s1"
$wsTranslations.Range("B4").Value = "s1"
$wsTranslations.Range("B5").Value = "s1"

$wsQuestion.Range("B2").Value = "s1"

# Give the freshly-inserted "Variable" column the same (manually sized) width
# as the Entity Id column next to it. The other, pre-existing columns keep
# whatever width Excel already auto-fit for them.
$wsTranslations.Columns.Item(2).ColumnWidth = 42.5
$wsQuestion.Columns.Item(2).ColumnWidth = 6.83

# Update selections/active sheet to match the final state of the workbook.
$wsTranslations.Activate()
$wsTranslations.Range("B6").Select()

$wsQuestion.Activate()
$wsQuestion.Range("B3").Select()

$wsTranslations.Activate()
